# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Add a new "ODI Bowling Extra" worksheet (after "ODI Batting Extra") with
#    MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.
# 2) Clean up placeholder empty cells that were left behind on the
#    "ODI Batting Extra" sheet (B/C/D/E on a handful of rows, plus E21).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Tidy up "ODI Batting Extra": drop the stray empty placeholder cells.
# ---------------------------------------------------------------------------
$wsBattingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyRanges = @("B5:E5", "B10:E10", "B14:E14", "B17:E17", "E21:E21")
foreach ($rng in $emptyRanges) {
    $wsBattingExtra.Range($rng).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet as the last tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBowlingExtra = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsBowlingExtra.Name = "ODI Bowling Extra"

# Reuse the bold/boxed header formatting already used on the other "Extra"
# sheet so the new header row matches the workbook's existing style.
$wsBattingExtra.Range("A1:C1").Copy()
$wsBowlingExtra.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats

$wsBowlingExtra.Range("A1").Value = "MATCH_CODE"
$wsBowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$wsBowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# row, MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL ($null => left blank)
$rows = @(
    @(2, "4093", "0", $null),
    @(3, "4097", "0", $null),
    @(4, "4106", "1", "20.00%"),
    @(5, "4107", "0", "40.00%"),
    @(6, "4111", $null, $null),
    @(7, "4113", "0", "10.00%"),
    @(8, "4143", "0", "20.00%"),
    @(9, "4150", $null, $null),
    @(10, "4156", "1", "10.00%"),
    @(11, "4158", "0", $null),
    @(12, "4190", "0", "10.00%"),
    @(13, "4284", "0", $null),
    @(14, "4285", "0", $null),
    @(15, "4299", "0", "10.00%"),
    @(16, "4301", "0", $null),
    @(17, "4347", "0", $null),
    @(18, "4352", $null, $null),
    @(19, "4442", "1", "10.00%"),
    @(20, "4444", "0", $null),
    @(21, "4446", "0", $null)
)

foreach ($item in $rows) {
    $r = $item[0]
    $matchCode = $item[1]
    $maidenOvers = $item[2]
    $percentWickets = $item[3]

    if ($null -ne $matchCode) {
        $cell = $wsBowlingExtra.Range("A" + $r)
        $cell.NumberFormat = "@"
        $cell.Value = $matchCode
    }
    if ($null -ne $maidenOvers) {
        $cell = $wsBowlingExtra.Range("B" + $r)
        $cell.NumberFormat = "@"
        $cell.Value = $maidenOvers
    }
    if ($null -ne $percentWickets) {
        $cell = $wsBowlingExtra.Range("C" + $r)
        $cell.NumberFormat = "@"
        $cell.Value = $percentWickets
    }
}

Write-Output "ODI Bowling Extra sheet added and ODI Batting Extra cleaned up"
